# Bid Altering (with suggestion and limit), Initial Balance set to 0, updated GUI
# (including time frame, indication when waiting)
#
# Adds a new results row (row 7) to the "testing of simple system" sheet,
# mirroring the existing rows but using a new bold style for the numeric
# cost/emission/time/distance columns and tagging the run with a new
# "<- Adjustable Bids" label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of results (row 7)
$ws.Range("A7").Value = "blon-1"
$ws.Range("B7").Value = 10
$ws.Range("C7").Value = 110
$ws.Range("D7").Value = 1850.04780660353
$ws.Range("E7").Value = 525.34671678699999
$ws.Range("F7").Value = 5105.9822323259996
$ws.Range("G7").Value = 713.36298178599998
$ws.Range("D7:G7").Font.Bold = $true
$ws.Range("H7").Value = 5
$ws.Range("I7").Value = "Car"
$ws.Range("J7").Value = "Emissions"
$ws.Range("K7").Value = "One Player"
$ws.Range("L7").Value = "<- Adjustable Bids"

# Move/update the active selection like in the edited workbook
[void]$ws.Range("B17").Select()
